# "Save Load I2C config implemented"
# Update the EEPROM memory-map sizing figures on the "Con I2C" (column D)
# and "Nuevo VCMC" (column C) sections of Hoja1 to reflect the new I2C
# config save/load data that now needs to be persisted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Column D ("Con I2C") byte-count updates ---
$ws.Range("D19").Value = 28                 # UserNames size grew by 1 byte
$ws.Range("D24").Value = 10                 # Gate config size
$ws.Range("D25").Value = 23                 # CV config size
$ws.Range("D26").Value = 18                 # Fader config size

# D21 ("Total Banco") is no longer derived from the formula; the new I2C
# save/load total is entered directly.
$ws.Range("D21").Value = 150

# --- Column C ("Nuevo VCMC") update ---
$ws.Range("C25").Value = 20                 # CV config size

# --- New row 38: free EEPROM bytes remaining once I2C config is saved ---
$ws.Range("D38").Formula = "=+D32-D34"

# Restore the sheet's visible selection to where the edit was made.
$ws.Range("D23").Select()
